# Trade #59 closed at 2026-02-16 21:33:18 - leadlag UP +0.000%
# Also reflects trade #31 (row 27 in "leadlag") closing out, and the new
# open trade #59 being logged (row 48 in "leadlag").

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force literal text so Excel doesn't auto-coerce percentages / dates / etc,
    # then restore the default "Normal" style so no stray number-format style
    # is left behind on the cell.
    if ($text -eq "") {
        # A bare quote prefix keeps the cell typed as Text with empty content
        # (plain "" would blank the cell out entirely / retype it as Number).
        $range.Value = "'"
        $range.Style = "Normal"
    } else {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.Style = "Normal"
    }
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Cells.Item(2,3).Value = 31
Set-TextValue $wsSummary.Cells.Item(2,4) "64.5%"
Set-TextValue $wsSummary.Cells.Item(2,5) "+6.6185%"
Set-TextValue $wsSummary.Cells.Item(2,6) "+0.2135%"

$wsSummary.Cells.Item(3,3).Value = 46
Set-TextValue $wsSummary.Cells.Item(3,4) "34.8%"
Set-TextValue $wsSummary.Cells.Item(3,5) "+4.8755%"
Set-TextValue $wsSummary.Cells.Item(3,6) "+0.1060%"

# ---------------------------------------------------------------------
# leadlag sheet
# ---------------------------------------------------------------------
$wsLead = $wb.Worksheets.Item("leadlag")

# Row 27 (Trade #31) - now closed
$wsLead.Cells.Item(27,7).Value = 69345.115123
Set-TextValue $wsLead.Cells.Item(27,8) "CLOSED"
$wsLead.Cells.Item(27,9).Value = 0.5794
$wsLead.Cells.Item(27,10).Value = 5.79
Set-TextValue $wsLead.Cells.Item(27,13) "time_exit_5min"
$wsLead.Cells.Item(27,14).Value = 5

# New row 48 (Trade #59) - newly opened
$wsLead.Cells.Item(48,1).Value = 59
Set-TextValue $wsLead.Cells.Item(48,2) "2026-02-16"
Set-TextValue $wsLead.Cells.Item(48,3) "21:33:18"
Set-TextValue $wsLead.Cells.Item(48,4) "leadlag"
Set-TextValue $wsLead.Cells.Item(48,5) "UP"
$wsLead.Cells.Item(48,6).Value = 68820.94500000001
Set-TextValue $wsLead.Cells.Item(48,7) ""
Set-TextValue $wsLead.Cells.Item(48,8) "OPEN"
$wsLead.Cells.Item(48,9).Value = 0
$wsLead.Cells.Item(48,10).Value = 0
$wsLead.Cells.Item(48,11).Value = 0.7358
Set-TextValue $wsLead.Cells.Item(48,12) "Coinbase leading with 0.074% move"
Set-TextValue $wsLead.Cells.Item(48,13) ""
$wsLead.Cells.Item(48,14).Value = 0

# ---------------------------------------------------------------------
# All Trades sheet - append closed Trade #31
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(32,1).Value = 31
Set-TextValue $wsAll.Cells.Item(32,2) "2026-02-16"
Set-TextValue $wsAll.Cells.Item(32,3) "21:28:05"
Set-TextValue $wsAll.Cells.Item(32,4) "leadlag"
Set-TextValue $wsAll.Cells.Item(32,5) "UP"
$wsAll.Cells.Item(32,6).Value = 68945.64
$wsAll.Cells.Item(32,7).Value = 69345.115123
Set-TextValue $wsAll.Cells.Item(32,8) "CLOSED"
$wsAll.Cells.Item(32,9).Value = 0.5794
$wsAll.Cells.Item(32,10).Value = 5.79
$wsAll.Cells.Item(32,11).Value = 0.75
Set-TextValue $wsAll.Cells.Item(32,12) "Coinbase leading with 0.088% move"
Set-TextValue $wsAll.Cells.Item(32,13) "time_exit_5min"
$wsAll.Cells.Item(32,14).Value = 5

# ---------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("Comparison")

$wsComp.Cells.Item(2,2).Value = 46
Set-TextValue $wsComp.Cells.Item(2,3) "34.8%"
Set-TextValue $wsComp.Cells.Item(2,4) "2.40"
Set-TextValue $wsComp.Cells.Item(2,5) "+0.5223%"
Set-TextValue $wsComp.Cells.Item(2,7) "1.50"
